$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 8.467854666666666
$ws.Range("H2").Value = 25.403564
$ws.Range("I2").Value = 0.1523462492674013
$ws.Range("J2").Value = 0.1523462492674013
$ws.Range("M2").Value = 13.76827833333333
$ws.Range("N2").Value = 41.304835
$ws.Range("O2").Value = 0.1794330173761795
$ws.Range("P2").Value = 0.1794330173761795
$ws.Range("Q2").Value = 116.5877799368822
$ws.Range("R2").Value = 1049.29001943194
$ws.Range("S2").Value = 0.02733594719199339
$ws.Range("T2").Value = 0.02733594719199339
$ws.Range("G3").Value = 8.467854666666666
$ws.Range("H3").Value = 25.403564
$ws.Range("I3").Value = 0.1523462492674013
$ws.Range("J3").Value = 0.1523462492674013
$ws.Range("O3").Value = 0.5028568908735593
$ws.Range("P3").Value = 0.5028568908735593
$ws.Range("Q3").Value = 326.7345630709675
$ws.Range("R3").Value = 2940.611067638707
$ws.Range("S3").Value = 0.07660836124285368
$ws.Range("T3").Value = 0.0766083612428537
$ws.Range("G4").Value = 8.467854666666666
$ws.Range("H4").Value = 25.403564
$ws.Range("I4").Value = 0.1523462492674013
$ws.Range("J4").Value = 0.1523462492674013
$ws.Range("O4").Value = 0.3177100917502611
$ws.Range("P4").Value = 0.3177100917502612
$ws.Range("Q4").Value = 206.4342159673422
$ws.Range("R4").Value = 1857.90794370608
$ws.Range("S4").Value = 0.04840194083255421
$ws.Range("T4").Value = 0.04840194083255422
$ws.Range("I5").Value = 0.1669927598427297
$ws.Range("J5").Value = 0.1669927598427297
$ws.Range("M5").Value = 13.76827833333333
$ws.Range("N5").Value = 41.304835
$ws.Range("O5").Value = 0.1794330173761795
$ws.Range("P5").Value = 0.1794330173761795
$ws.Range("Q5").Value = 127.7964848443617
$ws.Range("R5").Value = 1150.168363599255
$ws.Range("S5").Value = 0.02996401477855669
$ws.Range("T5").Value = 0.0299640147785567
$ws.Range("I6").Value = 0.1669927598427297
$ws.Range("J6").Value = 0.1669927598427297
$ws.Range("O6").Value = 0.5028568908735593
$ws.Range("P6").Value = 0.5028568908735593
$ws.Range("S6").Value = 0.08397346001291005
$ws.Range("T6").Value = 0.08397346001291005
$ws.Range("I7").Value = 0.1669927598427297
$ws.Range("J7").Value = 0.1669927598427297
$ws.Range("O7").Value = 0.3177100917502611
$ws.Range("P7").Value = 0.3177100917502612
$ws.Range("S7").Value = 0.05305528505126298
$ws.Range("T7").Value = 0.05305528505126299
$ws.Range("I8").Value = 0.680660990889869
$ws.Range("J8").Value = 0.680660990889869
$ws.Range("M8").Value = 13.76827833333333
$ws.Range("N8").Value = 41.304835
$ws.Range("O8").Value = 0.1794330173761795
$ws.Range("P8").Value = 0.1794330173761795
$ws.Range("Q8").Value = 520.8973256584717
$ws.Range("R8").Value = 4688.075930926245
$ws.Range("S8").Value = 0.1221330554056294
$ws.Range("T8").Value = 0.1221330554056294
$ws.Range("I9").Value = 0.680660990889869
$ws.Range("J9").Value = 0.680660990889869
$ws.Range("O9").Value = 0.5028568908735593
$ws.Range("P9").Value = 0.5028568908735593
$ws.Range("S9").Value = 0.3422750696177956
$ws.Range("T9").Value = 0.3422750696177956
$ws.Range("I10").Value = 0.680660990889869
$ws.Range("J10").Value = 0.680660990889869
$ws.Range("O10").Value = 0.3177100917502611
$ws.Range("P10").Value = 0.3177100917502612
$ws.Range("S10").Value = 0.2162528658664439
$ws.Range("T10").Value = 0.216252865866444
